$p = $ppt.ActivePresentation

# Add a new slide (slide 2) after the existing slide, using the
# "Title and Content" layout (CustomLayouts index 2 on the slide master).
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$s2 = $p.Slides.AddSlide(2, $layout)

# The new slide gets a Title placeholder and a Content placeholder,
# both left empty (a blank "started" slide).
$title = $s2.Shapes.Item(1)
$content = $s2.Shapes.Item(2)

$content.Name = "Content Placeholder 10"

$title.TextFrame.TextRange.LanguageID = "en-US"
$content.TextFrame.TextRange.LanguageID = "en-US"
